# Add two new columns (I and J) with header labels "I0" and "IF",
# matching the formatting of the existing header cells, plus their
# corresponding data values (9 and 9) in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, alignment) from the last
# existing header cell (H1) onto the two new header cells so they
# look consistent with the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header text.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set the new data values for row 2.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
